# Implement writeemail, bugs identified in README.md
# Adds three new "Postman" test rows (9, 10, 11) to the email log sheet,
# mirroring the records captured while exercising the writeemail endpoint.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (row 10) down onto the
# three new rows so the bold/bordered ID column and the date-time number
# format on the "Date Sent" column carry through.
$ws.Range("A10:I10").Copy()
$ws.Range("A11:I13").PasteSpecial(-4122)

# Row 11 - id 9 / "Postman4"
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = "Postman4"
$ws.Cells.Item(11, 3).Value = 25569.019386574073
$ws.Cells.Item(11, 4).Value = "AGAIN I sent yet ANOTHER email via postman"
$ws.Cells.Item(11, 5).Value = "marleevaughn@outlook.com"
$ws.Cells.Item(11, 6).Value = "Marlee Vaughn"
$ws.Cells.Item(11, 7).Value = "duanevaughn@hotmail.com"
$ws.Cells.Item(11, 8).Value = "Duane Vaughn"
$ws.Cells.Item(11, 9).Value = $true

# Row 12 - id 10 / "Postman5" (this one saved as a draft, Draft = FALSE)
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = "Postman5"
$ws.Cells.Item(12, 3).Value = 25569.019386574073
$ws.Cells.Item(12, 4).Value = "updated AGAIN I sent yet ANOTHER email via postman"
$ws.Cells.Item(12, 5).Value = "marleevaughn@outlook.com"
$ws.Cells.Item(12, 6).Value = "Marlee Vaughn"
$ws.Cells.Item(12, 7).Value = "duanevaughn@hotmail.com"
$ws.Cells.Item(12, 8).Value = "Duane Vaughn"
$ws.Cells.Item(12, 9).Value = $false

# Row 13 - id 11 / "Postman6"
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = "Postman6"
$ws.Cells.Item(13, 3).Value = 25569.019389594909
$ws.Cells.Item(13, 4).Value = "AGAIN I sent yet ANOTHER email via postman"
$ws.Cells.Item(13, 5).Value = "marleevaughn@outlook.com"
$ws.Cells.Item(13, 6).Value = "Marlee Vaughn"
$ws.Cells.Item(13, 7).Value = "duanevaughn@hotmail.com"
$ws.Cells.Item(13, 8).Value = "Duane Vaughn"
$ws.Cells.Item(13, 9).Value = $true

# Re-fit column C ("Date Sent") now that it holds more rows, and drop the
# now-unneeded custom width on column D.
$ws.Columns.Item(3).AutoFit()
$ws.Columns.Item(4).ColumnWidth = $ws.StandardWidth

# Move the active selection to reflect where editing left off.
$ws.Range("C11").Select() | Out-Null
